# Update weekly fruit/vegetable price records: rotate the data among
# rows 3, 4 and 5 (date, volume, min/max/weighted price, price per Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44253
$ws.Range("M3").Value = 160

# Row 4
$ws.Range("D4").Value = 44252
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 13000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 13500
$ws.Range("S4").Value = 750

# Row 5
$ws.Range("D5").Value = 44257
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("S5").Value = 806
